$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# The answer grid occupies rows 1, 5, 9, 13, 17 (each followed by 3 blank
# rows reserved for handwriting) x columns 1-5. Cell text is replaced by
# position rather than Find/Replace because several old values repeat
# (e.g. "12÷8=1, 4" and "40÷3=13, 1" each appear twice).

# Row 1, Col 1: "13÷8=1, 5" -> "76÷6=12, 4"
$t.Cell(1,1).Range.Text = "76÷6=12, 4"
# Row 1, Col 2: "53÷2=26, 1" -> "15÷8=1, 7"
$t.Cell(1,2).Range.Text = "15÷8=1, 7"
# Row 1, Col 3: "69÷3=23, 0" -> "20÷8=2, 4"
$t.Cell(1,3).Range.Text = "20÷8=2, 4"
# Row 1, Col 4: "51÷8=6, 3" -> "82÷4=20, 2"
$t.Cell(1,4).Range.Text = "82÷4=20, 2"
# Row 1, Col 5: "22÷3=7, 1" -> "52÷4=13, 0"
$t.Cell(1,5).Range.Text = "52÷4=13, 0"

# Row 5, Col 1: "18÷2=9, 0" -> "61÷2=30, 1"
$t.Cell(5,1).Range.Text = "61÷2=30, 1"
# Row 5, Col 2: "16÷4=4, 0" -> "29÷9=3, 2"
$t.Cell(5,2).Range.Text = "29÷9=3, 2"
# Row 5, Col 3: "12÷8=1, 4" -> "94÷7=13, 3"
$t.Cell(5,3).Range.Text = "94÷7=13, 3"
# Row 5, Col 4: "68÷6=11, 2" -> "45÷6=7, 3"
$t.Cell(5,4).Range.Text = "45÷6=7, 3"
# Row 5, Col 5: "92÷9=10, 2" -> "43÷3=14, 1"
$t.Cell(5,5).Range.Text = "43÷3=14, 1"

# Row 9, Col 1: "44÷2=22, 0" -> "88÷5=17, 3"
$t.Cell(9,1).Range.Text = "88÷5=17, 3"
# Row 9, Col 2: "40÷3=13, 1" -> "47÷2=23, 1"
$t.Cell(9,2).Range.Text = "47÷2=23, 1"
# Row 9, Col 3: "61÷8=7, 5" -> "87÷3=29, 0"
$t.Cell(9,3).Range.Text = "87÷3=29, 0"
# Row 9, Col 4: "76÷9=8, 4" -> "13÷9=1, 4"
$t.Cell(9,4).Range.Text = "13÷9=1, 4"
# Row 9, Col 5: "57÷3=19, 0" -> "41÷9=4, 5"
$t.Cell(9,5).Range.Text = "41÷9=4, 5"

# Row 13, Col 1: "18÷3=6, 0" -> "25÷3=8, 1"
$t.Cell(13,1).Range.Text = "25÷3=8, 1"
# Row 13, Col 2: "63÷2=31, 1" -> "39÷8=4, 7"
$t.Cell(13,2).Range.Text = "39÷8=4, 7"
# Row 13, Col 3: "24÷5=4, 4" -> "46÷8=5, 6"
$t.Cell(13,3).Range.Text = "46÷8=5, 6"
# Row 13, Col 4: "65÷5=13, 0" -> "27÷8=3, 3"
$t.Cell(13,4).Range.Text = "27÷8=3, 3"
# Row 13, Col 5: "12÷8=1, 4" -> "56÷5=11, 1"
$t.Cell(13,5).Range.Text = "56÷5=11, 1"

# Row 17, Col 1: "54÷5=10, 4" -> "51÷5=10, 1"
$t.Cell(17,1).Range.Text = "51÷5=10, 1"
# Row 17, Col 2: "55÷2=27, 1" -> "72÷7=10, 2"
$t.Cell(17,2).Range.Text = "72÷7=10, 2"
# Row 17, Col 3: "40÷3=13, 1" -> "76÷3=25, 1"
$t.Cell(17,3).Range.Text = "76÷3=25, 1"
# Row 17, Col 4: "36÷9=4, 0" -> "17÷9=1, 8"
$t.Cell(17,4).Range.Text = "17÷9=1, 8"
# Row 17, Col 5: "64÷5=12, 4" -> "55÷7=7, 6"
$t.Cell(17,5).Range.Text = "55÷7=7, 6"

